$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The document currently ends with a trailing empty paragraph that
# only holds the "_GoBack" bookmark. We need to:
#   1) insert a brand-new paragraph ("One of the topics...") right
#      before that bookmark paragraph, and
#   2) turn the bookmark paragraph itself into the "Another problem
#      at work..." paragraph, with four runs of text BEFORE the
#      bookmark and one run AFTER it.
#
# The backend silently coalesces any two adjacent runs that share the
# same (empty/default) formatting, so simply typing the four "before"
# runs straight into the same paragraph would merge them into a
# single <w:r>, which would not match the source diff. To keep them
# distinct, each chunk is built in its own temporary paragraph
# (created with InsertParagraphBefore, which leaves the bookmark
# attached to the original paragraph) and the paragraph marks that
# separated the scratch paragraphs are deleted afterwards. That join
# does NOT trigger the same-format run merge, so the run boundaries
# created while authoring survive into the saved document - matching
# how Word actually produced this file (distinct typing/edit
# sessions yield distinct runs).
# ------------------------------------------------------------------

$lastIndex = $d.Paragraphs.Count
$bookmarkPara = $d.Paragraphs.Item($lastIndex)

# --- Step 1: brand new standalone paragraph before the bookmark para ---
$introText = 'One of the topics on deep learning I would like to explore are GANs (Generative Adversarial Networks). The main reason for this is that at work (Boeing) some teams are trying to implement this to obtain a more accurate 3D models of the airplanes that we are building by taking some of the pictures (2D images) that we have and generate a 3D more realistic depiction of the airplane. Another related project that is currently under develop is the generation a realistic 3D model of the different landing runways, this with the purpose of improve landing simulations for pilots training.'

$bookmarkPara.Range.InsertParagraphBefore()
$introPara = $d.Paragraphs.Item($lastIndex)
$introPara.Range.Text = $introText

# $lastIndex now points at the (still-empty) bookmark paragraph again,
# one slot further down.
$lastIndex = $lastIndex + 1
$bookmarkPara = $d.Paragraphs.Item($lastIndex)

# --- Step 2: four runs before the bookmark, one after it ---
$run1 = 'Another problem at work '
$run2 = 'where'
$run3 = ' I think Deep Learning could help with is the identification of damage on different parts of the airplane. By now we have a substantial library of pictures and descriptions of different kind of damage on different parts of the '
$run4 = 'airplane. This is a very time consuming task and based on the current database of pictures of good parts (per drawing parts) and damaged parts (with labels) we could potentially streamline the process by using deep learning to first determine if there is damage on a part by looking at it and then to identify the k'
$run5 = 'ind of damage that is seen on that particular part. '

# Create four empty scratch paragraphs immediately before the bookmark
# paragraph (InsertParagraphBefore keeps the bookmark on the original
# paragraph, shifting it one slot further down each time).
for ($i = 0; $i -lt 4; $i++) {
    $bookmarkPara.Range.InsertParagraphBefore()
}

$p1 = $d.Paragraphs.Item($lastIndex)
$p1.Range.Text = $run1
$p2 = $d.Paragraphs.Item($lastIndex + 1)
$p2.Range.Text = $run2
$p3 = $d.Paragraphs.Item($lastIndex + 2)
$p3.Range.Text = $run3
$p4 = $d.Paragraphs.Item($lastIndex + 3)
$p4.Range.Text = $run4

# The bookmark paragraph is now the 5th of this block; append the
# trailing run AFTER the bookmark, inside that same paragraph.
$p5 = $d.Paragraphs.Item($lastIndex + 4)
$p5.Range.InsertAfter($run5)

# --- Step 3: join the five scratch paragraphs back into a single ---
# --- paragraph by deleting the paragraph marks between them.     ---
for ($i = 0; $i -lt 4; $i++) {
    $p = $d.Paragraphs.Item($lastIndex)
    $mark = $d.Range($p.Range.End - 1, $p.Range.End)
    $mark.Delete()
}
